$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.581.47'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('D3').Value = '1.868.17'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').Value = '''324.72'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').Value = '''0.4615'
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('D8').Value = '''0.3878'
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '''0.07860'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '''0.9740'
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('D11').Value = '''21.92'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '1.885.42'
$ws.Range('E12').Value = '  +2.71%  '
$ws.Range('D13').Value = '''6.982'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').Value = '''5.693'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').Value = '''0.06957'
$ws.Range('E15').Value = '  +2.54%  '
$ws.Range('D16').Value = '''88.05'
$ws.Range('E16').Value = '  +0.92%  '
$ws.Range('D17').Value = '''1.005'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').Value = '''0.00001001'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').Value = '''16.80'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').Value = '''1.003'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = '28.598.30'
$ws.Range('E21').Value = '  +2.05%  '
$ws.Range('D22').Value = '''5.272'
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = '''10.99'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').Value = '2.105.58'
$ws.Range('E25').Value = '  +2.35%  '
$ws.Range('D26').Value = '''152.50'
$ws.Range('E26').Value = '  -0.85%  '
$ws.Range('D27').Value = '''19.21'
$ws.Range('E27').Value = '  +0.28%  '
$ws.Range('D28').Value = '''5.772'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').Value = '''1.984'
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('D30').Value = '''119.19'
$ws.Range('E30').Value = '  +1.51%  '
$ws.Range('D31').Value = '''0.09336'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').Value = '''0.9143'
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').Value = '''5.258'
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').Value = '''1.335'
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range('D35').Value = '''3.326'
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').Value = '''0.05776'
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02103'
$ws.Range('E37').Value = '  -1.81%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.155'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').Value = '''7.722'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').Value = '''0.5624'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('D41').Value = '''0.1783'
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('D42').Value = '''9.763'
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('D43').Value = '''0.07179'
$ws.Range('E43').Value = '  +2.53%  '
$ws.Range('D44').Value = '''11.76'
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('D45').Value = '''0.5297'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('D46').Value = '''2.149'
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('D47').Value = '''1.142'
$ws.Range('E47').Value = '  +2.61%  '
$ws.Range('D48').Value = '''1.828'
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').Value = '''112.90'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '''2.406'
$ws.Range('E50').Value = '  +3.72%  '
$ws.Range('E51').Value = '  +0.28%  '
